$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "position" column (B) previously held the 1-based row/position index
# for each attribute (1..20). Data processing was cleaned up so that the
# position is no longer computed per-row here; instead every row reports
# -1, signalling "position not set / computed elsewhere".
$ws.Range("B2:B21").Value = -1
